$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data refresh (new day's figures) ---
# A2: DateVal            43928 -> 43929
# B2: TotalUKCases        55242 -> 60733
# C2: NewUKCases            3634 -> 5491
# D2: TotalUKDeaths        6159 -> 7097
# E2: DailyUKDeaths          786 -> 938
# F2: EnglandCases        42990 -> 50756
# G2: EnglandDeaths         4897 -> 6483
# H2: ScotlandCases (fix)   3961 -> 4565
# I2: ScotlandDeaths (fix)   220 -> 296
# J2: WalesCases            3499 -> 4073
# K2: WalesDeaths             193 -> 245
# L2: NICases                1158 -> 1339
# M2: NIDeaths                  63 -> 73

$ws.Range("A2").Value = 43929
$ws.Range("B2").Value = 60733
$ws.Range("C2").Value = 5491
$ws.Range("D2").Value = 7097
$ws.Range("E2").Value = 938
$ws.Range("F2").Value = 50756
$ws.Range("G2").Value = 6483
$ws.Range("H2").Value = 4565
$ws.Range("I2").Value = 296
$ws.Range("J2").Value = 4073
$ws.Range("K2").Value = 245
$ws.Range("L2").Value = 1339
$ws.Range("M2").Value = 73

# Keep the custom "##,##0" integer display format on the data row.
$ws.Range("B2:M2").NumberFormat = "##,##0"

# Move the active selection to reflect where editing left off.
$ws.Range("C9").Select() | Out-Null
